# Auto-generated Excel COM-interop script
# Applies the data updates described by the upstream diff to Sheets/Bahamut_Profits.xlsx
# (scheduled runner refresh of currentAveragePrice / Leve profit columns).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 262.05
$ws.Range("I33").Value = 143.33333
$ws.Range("J33").Value = 312.92856
$ws.Range("K33").Value = 143.33333
$ws.Range("L33").Value = 312.92856
$ws.Range("M33").Value = 85.66667000000001
$ws.Range("N33").Value = -770.9285600000001
$ws.Range("H40").Value = 2041.6666
$ws.Range("J40").Value = 2041.6666
$ws.Range("L40").Value = 2041.6666
$ws.Range("N40").Value = -2391.6666
$ws.Range("H64").Value = 5259.8
$ws.Range("J64").Value = 4800
$ws.Range("L64").Value = 4800
$ws.Range("N64").Value = -5296
$ws.Range("H67").Value = 5259.8
$ws.Range("J67").Value = 4800
$ws.Range("L67").Value = 4800
$ws.Range("N67").Value = -6516
$ws.Range("H74").Value = 3247907.8
$ws.Range("I74").Value = 3996532.2
$ws.Range("K74").Value = 3996532.2
$ws.Range("M74").Value = -3995596.2
$ws.Range("H77").Value = 3247907.8
$ws.Range("I77").Value = 3996532.2
$ws.Range("K77").Value = 19982661
$ws.Range("M77").Value = -19977981
$ws.Range("H103").Value = 747.4737
$ws.Range("I103").Value = 816.5333000000001
$ws.Range("J103").Value = 488.5
$ws.Range("K103").Value = 2449.5999
$ws.Range("L103").Value = 1465.5
$ws.Range("M103").Value = -1863.5999
$ws.Range("N103").Value = -2637.5
$ws.Range("H126").Value = 37000
$ws.Range("J126").Value = 37000
$ws.Range("L126").Value = 37000
$ws.Range("N126").Value = -46880
$ws.Range("H132").Value = 1819971.6
$ws.Range("I132").Value = 1565.5366
$ws.Range("J132").Value = 7145303.5
$ws.Range("K132").Value = 4696.6098
$ws.Range("L132").Value = 21435910.5
$ws.Range("M132").Value = -2166.6098
$ws.Range("N132").Value = -21440970.5
$ws.Range("H135").Value = 1187.4814
$ws.Range("I135").Value = 1189.4615
$ws.Range("J135").Value = 1136
$ws.Range("K135").Value = 10705.1535
$ws.Range("L135").Value = 10224
$ws.Range("M135").Value = -8170.153499999999
$ws.Range("N135").Value = -15294
$ws.Range("H137").Value = 996.46155
$ws.Range("I137").Value = 900.55554
$ws.Range("K137").Value = 2701.66662
$ws.Range("M137").Value = -151.66662
$ws.Range("H138").Value = 3062.8367
$ws.Range("I138").Value = 1334.1936
$ws.Range("J138").Value = 3862.6567
$ws.Range("K138").Value = 4002.5808
$ws.Range("L138").Value = 11587.9701
$ws.Range("M138").Value = 1137.4192
$ws.Range("N138").Value = -21867.9701
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3938.99
$ws.Range("I32").Value = 3769.7812
$ws.Range("J32").Value = 8000
$ws.Range("K32").Value = 3769.7812
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = -3482.7812
$ws.Range("N32").Value = -8574
$ws.Range("H61").Value = 2646.7917
$ws.Range("I61").Value = 2635.95
$ws.Range("J61").Value = 2701
$ws.Range("K61").Value = 2635.95
$ws.Range("L61").Value = 2701
$ws.Range("M61").Value = -2423.95
$ws.Range("N61").Value = -3125
$ws.Range("H122").Value = 890
$ws.Range("I122").Value = 890
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2670
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -220
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1785.7441
$ws.Range("I132").Value = 1308.4849
$ws.Range("J132").Value = 3360.7
$ws.Range("K132").Value = 3925.4547
$ws.Range("L132").Value = 10082.1
$ws.Range("M132").Value = -1395.4547
$ws.Range("N132").Value = -15142.1
$ws.Range("H136").Value = 2646.7917
$ws.Range("I136").Value = 2635.95
$ws.Range("J136").Value = 2701
$ws.Range("K136").Value = 7907.849999999999
$ws.Range("L136").Value = 8103
$ws.Range("M136").Value = -5357.849999999999
$ws.Range("N136").Value = -13203
$ws.Range("H139").Value = 55571.668
$ws.Range("J139").Value = 55571.668
$ws.Range("L139").Value = 55571.668
$ws.Range("N139").Value = -65851.66800000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 60000
$ws.Range("J59").Value = 60000
$ws.Range("L59").Value = 60000
$ws.Range("N59").Value = -61694
$ws.Range("H86").Value = 66677400
$ws.Range("I86").Value = 100013010
$ws.Range("J86").Value = 6196
$ws.Range("K86").Value = 100013010
$ws.Range("L86").Value = 6196
$ws.Range("M86").Value = -100011887
$ws.Range("N86").Value = -8442
$ws.Range("H89").Value = 66677400
$ws.Range("I89").Value = 100013010
$ws.Range("J89").Value = 6196
$ws.Range("K89").Value = 500065050
$ws.Range("L89").Value = 30980
$ws.Range("M89").Value = -500059434
$ws.Range("N89").Value = -42212
$ws.Range("H107").Value = 10976.077
$ws.Range("I107").Value = 824.0833
$ws.Range("K107").Value = 824.0833
$ws.Range("M107").Value = 1095.9167
$ws.Range("H134").Value = 1968.6227
$ws.Range("I134").Value = 1573.5682
$ws.Range("J134").Value = 3900
$ws.Range("K134").Value = 4720.7046
$ws.Range("L134").Value = 11700
$ws.Range("M134").Value = -2185.7046
$ws.Range("N134").Value = -16770
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 68027.56
$ws.Range("I31").Value = 7672.5713
$ws.Range("J31").Value = 114970.336
$ws.Range("K31").Value = 7672.5713
$ws.Range("L31").Value = 114970.336
$ws.Range("N31").Value = -115560.336
$ws.Range("M31").Value = -7377.5713
$ws.Range("H34").Value = 68027.56
$ws.Range("I34").Value = 7672.5713
$ws.Range("J34").Value = 114970.336
$ws.Range("K34").Value = 7672.5713
$ws.Range("L34").Value = 114970.336
$ws.Range("N34").Value = -115374.336
$ws.Range("M34").Value = -7470.5713
$ws.Range("H122").Value = 1159.2
$ws.Range("I122").Value = 1018.6667
$ws.Range("K122").Value = 3056.0001
$ws.Range("M122").Value = -606.0001000000002
$ws.Range("H132").Value = 973.3958
$ws.Range("I132").Value = 803.04877
$ws.Range("J132").Value = 1971.1428
$ws.Range("K132").Value = 2409.14631
$ws.Range("L132").Value = 5913.428400000001
$ws.Range("M132").Value = 120.8536899999999
$ws.Range("N132").Value = -10973.4284
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3921748
$ws.Range("J2").Value = 58.4
$ws.Range("L2").Value = 350.4
$ws.Range("N2").Value = -576.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4530.727
$ws.Range("I70").Value = 4259.8
$ws.Range("J70").Value = 4756.5
$ws.Range("K70").Value = 4259.8
$ws.Range("L70").Value = 4756.5
$ws.Range("M70").Value = -3989.8
$ws.Range("N70").Value = -5296.5
$ws.Range("H73").Value = 4530.727
$ws.Range("I73").Value = 4259.8
$ws.Range("J73").Value = 4756.5
$ws.Range("K73").Value = 4259.8
$ws.Range("L73").Value = 4756.5
$ws.Range("M73").Value = -3323.8
$ws.Range("N73").Value = -6628.5
$ws.Range("H107").Value = 293.55173
$ws.Range("I107").Value = 187.77777
$ws.Range("J107").Value = 466.63635
$ws.Range("K107").Value = 187.77777
$ws.Range("L107").Value = 466.63635
$ws.Range("M107").Value = 1732.22223
$ws.Range("N107").Value = -4306.63635
$ws.Range("H132").Value = 2220.5283
$ws.Range("I132").Value = 1835.8536
$ws.Range("J132").Value = 3534.8333
$ws.Range("K132").Value = 5507.560799999999
$ws.Range("L132").Value = 10604.4999
$ws.Range("M132").Value = -2977.560799999999
$ws.Range("N132").Value = -15664.4999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 570.64703
$ws.Range("I22").Value = 337.375
$ws.Range("J22").Value = 778
$ws.Range("K22").Value = 337.375
$ws.Range("L22").Value = 778
$ws.Range("M22").Value = -42.375
$ws.Range("N22").Value = -1368
$ws.Range("H27").Value = 570.64703
$ws.Range("I27").Value = 337.375
$ws.Range("J27").Value = 778
$ws.Range("K27").Value = 337.375
$ws.Range("L27").Value = 778
$ws.Range("M27").Value = -230.375
$ws.Range("N27").Value = -992
$ws.Range("H46").Value = 28572354
$ws.Range("I46").Value = 40001096
$ws.Range("J46").Value = 492
$ws.Range("K46").Value = 40001096
$ws.Range("L46").Value = 492
$ws.Range("M46").Value = -40000908
$ws.Range("N46").Value = -868
$ws.Range("H132").Value = 2555.75
$ws.Range("I132").Value = 1845.1072
$ws.Range("J132").Value = 4545.55
$ws.Range("K132").Value = 5535.321599999999
$ws.Range("L132").Value = 13636.65
$ws.Range("M132").Value = -3005.321599999999
$ws.Range("N132").Value = -18696.65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1121.5538
$ws.Range("I132").Value = 901.2558
$ws.Range("J132").Value = 1552.1364
$ws.Range("K132").Value = 2703.7674
$ws.Range("L132").Value = 4656.4092
$ws.Range("M132").Value = -173.7674000000002
$ws.Range("N132").Value = -9716.4092
